# Rostock Printer BOM update: adding stl's and updated xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rostock Max V2")

# --- Row 6: new line item "12"x12" Ultem Sheet (PEI)" / "McMaster-Carr, for Boro Glass" ---
$ws.Range("A6").Value = "12""x12"" Ultem Sheet (PEI)"
$ws.Range("B6").Value = "McMaster-Carr, for Boro Glass"
$ws.Range("C6").Value = 19
$ws.Range("D6").Value = 1

# --- Row 21: split the hotend line into a short item name + separate description ---
$ws.Range("B21").Value = "Fully Assembled 1.75mm Universal (with Bowden add-on) (24v)"
$ws.Range("A21").Value = "E3D All-metal v6 HotEnd"

# --- Indent the sub-item labels under the "Duet Upgrade" section (rows 15-19) ---
$ws.Cells.Item(15, 1).HorizontalAlignment = -4131
$ws.Cells.Item(15, 1).IndentLevel = 1
$ws.Range("A15").Copy()
$ws.Range("A16:A19").PasteSpecial(-4122)

# --- Column A widens slightly to fit the new, longer item text ---
$ws.Columns("A:A").ColumnWidth = 21.425

# --- Restore the cursor/selection to where it was left in the saved file ---
$ws.Range("A25").Select()
